# Moving.docx - "Added Alternate Courses to Moving"
#
# The existing "Alternate Courses" paragraph ("Movement is against wall: The
# character does not move. ") is rewritten into a richer "Character is
# against a wall or enemy: ..." alternate course, and a brand-new "Character
# is frozen or paralyzed: The character does not move. " alternate course is
# appended right after it.
#
# Because the new wording mixes bold/non-bold runs, re-introduces the
# _GoBack bookmark in the middle of the sentence, and splits "the" across two
# runs ("t" + "he") exactly as the canonical OOXML does, we rebuild the
# paragraph (and its successor) from a WordprocessingML fragment via
# Range.InsertXML rather than a sequence of Find/Replace calls - that keeps
# run boundaries, bookmark placement and the paragraph-mark formatting
# (bold + sz 24) bit-for-bit faithful to the target.

$d = $word.ActiveDocument

# Locate the "Alternate Courses:" content paragraph - the one that currently
# reads "Movement is against wall: The character does not move. "
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Movement is against wall:")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Movement is against wall:' paragraph"
}

$fragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
  '<w:pPr><w:rPr><w:b/><w:sz w:val="24"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Character is against a wall</w:t></w:r>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> or </w:t></w:r>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>enemy</w:t></w:r>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' +
  '<w:r><w:t>If the movement is in the direction of the wall</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> or </w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:t>enemy</w:t></w:r>' +
  '<w:r><w:t>, t</w:t></w:r>' +
  '<w:r><w:t>he character does not move.</w:t></w:r>' +
  '<w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
'</w:p>' +
'<w:p>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Character is frozen or paralyzed: </w:t></w:r>' +
  '<w:r><w:t>The character does not move.</w:t></w:r>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

$target.Range.InsertXML($fragment)

# NOTE: the upstream diff also drops a bare <w:semiHidden/> flag from the
# built-in "Default Paragraph Font" character style in styles.xml. Word's
# object model (real Word included) does not expose a SemiHidden property
# on the Style object - Style.Visibility/Hidden maps to <w:hidden/>, a
# different flag - and Document/Range.WordOpenXML is read-only here, so
# that particular style-table byte is outside what Range.InsertXML /
# Style.* can reach from script and is intentionally left untouched.
